$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (the first data row), shifting existing data down by one row.
$ws.Rows.Item(2).Insert()
$ws.Range("A2:R2").ClearFormats()

# Populate the new row 2 with the new weekly data. Columns A,B,C,E,F,G,H,I,N,O,Q,R are
# constant across the whole sheet, so copy them down from the row that is now row 3
# (the former row 2), and set the varying columns (D,J,K,L,M,P) to the new values.
$ws.Cells.Item(2, 1).Value2 = $ws.Cells.Item(3, 1).Value2
$ws.Cells.Item(2, 2).Value2 = $ws.Cells.Item(3, 2).Value2
$ws.Cells.Item(2, 3).Value2 = $ws.Cells.Item(3, 3).Value2
$ws.Cells.Item(2, 4).Value2 = 45160
$ws.Cells.Item(2, 4).NumberFormat = $ws.Cells.Item(3, 4).NumberFormat
$ws.Cells.Item(2, 5).Value2 = $ws.Cells.Item(3, 5).Value2
$ws.Cells.Item(2, 6).Value2 = $ws.Cells.Item(3, 6).Value2
$ws.Cells.Item(2, 7).Value2 = $ws.Cells.Item(3, 7).Value2
$ws.Cells.Item(2, 8).Value2 = $ws.Cells.Item(3, 8).Value2
$ws.Cells.Item(2, 9).Value2 = $ws.Cells.Item(3, 9).Value2
$ws.Cells.Item(2, 10).Value2 = 35
$ws.Cells.Item(2, 11).Value2 = 13000
$ws.Cells.Item(2, 12).Value2 = 13000
$ws.Cells.Item(2, 13).Value2 = 13000
$ws.Cells.Item(2, 14).Value2 = $ws.Cells.Item(3, 14).Value2
$ws.Cells.Item(2, 15).Value2 = $ws.Cells.Item(3, 15).Value2
$ws.Cells.Item(2, 16).Value2 = 1300
$ws.Cells.Item(2, 17).Value2 = $ws.Cells.Item(3, 17).Value2
$ws.Cells.Item(2, 18).Value2 = $ws.Cells.Item(3, 18).Value2
